# Add a "Formula Text" documentation column (H) that shows the literal
# formula text of the PV(...) examples in column F, using FORMULATEXT().

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column, using the same bold/centered formatting as
# the other header cells in row 1 (same style as A1:F1's "Present Value" etc.)
$ws.Range("H1").Value = "Formula Text"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108

# One FORMULATEXT() formula per PV example row, mirroring the column F
# formula on the same row.
$ws.Range("H2").Formula = "=FORMULATEXT(F2)"
$ws.Range("H3").Formula = "=FORMULATEXT(F3)"
$ws.Range("H4").Formula = "=FORMULATEXT(F4)"
$ws.Range("H5").Formula = "=FORMULATEXT(F5)"
$ws.Range("H6").Formula = "=FORMULATEXT(F6)"

# Widen the new column so the formula text is fully visible.
$ws.Columns.Item(8).ColumnWidth = 30.71

# Move the active selection to F2 (matches the resulting workbook state).
[void]$ws.Range("F2").Select()
